$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (pushes existing rows 2..76 down to 3..77)
$ws.Rows("2:2").Insert()

# The newly inserted row inherits the header row's bold/border style by default.
# Reset it to a plain (unstyled) row like the other data rows, then restore the
# date number format on column D to match the rest of the date column.
$ws.Rows("2:2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the new data record
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44496
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100102
$ws.Range("H2").Value = "Cítricos"
$ws.Range("I2").Value = 100102004
$ws.Range("J2").Value = "Mandarina"
$ws.Range("K2").Value = "Murcott"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("Q2").Value = '$/caja 20 kilos'
$ws.Range("R2").Value = "Región de Coquimbo"
$ws.Range("S2").Value = 675
$ws.Range("T2").Value = 20
